$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C5) from 2023-10-22 (45221) to 2023-10-25 (45224)
$ws.Range("C2:C5").Value = 45224
